$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 8; $r++) {
    $ws.Cells.Item($r, 1).Value = 3
    $ws.Cells.Item($r, 2).Value = 2024
}

$ws.Range("H10").Select() | Out-Null
